$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting existing rows 3-6 down to 4-7
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with data (copy constant columns from row 2 pattern, set new values)
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C3").Value = "Los Lagos"
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("D3").Value = 44498
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = "Chirimoya"
$ws.Range("K3").Value = "Cultivar IV Región"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("Q3").Value = "$/bandeja 8 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 2438
$ws.Range("T3").Value = 8
